# Atualizando o arquivo XLSX
# Apply updates to row 2 odds values on Sheet1

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("J2").Value = 2.5
$ws.Range("W2").Value = 11
$ws.Range("AH2").Value = 15
$ws.Range("AO2").Value = 10
